# Update cryptos "Price" (D) and "Volume(1h)" (E) columns per Feb 23 2023 GitHub Actions data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.908.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.647.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3887"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.25%  "
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.19"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.346"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.002"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08414"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.051"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.851"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001314"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.648.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06978"
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.910"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.904.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.448"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.973"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.389"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "138.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.775"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.831.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.043"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08026"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02947"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.697"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2674"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09084"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7541"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  -2.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.417"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6918"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.450"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.077"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08273"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.204"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.24%  "
